$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Weekly Quantity" ---
$ws1 = $wb.Worksheets.Item(1)

# Delete rows from the bottom up so earlier row numbers stay stable while we work:
#   row 26 -> date 45375.99999999999 (Order Week), qty 10
#   row 14 -> date 45137.99999999999, qty 40
#   row 11 -> date 45116.99999999999, qty 300
$ws1.Rows.Item(26).Delete()
$ws1.Rows.Item(14).Delete()
$ws1.Rows.Item(11).Delete()

# After the deletions above, the row that used to be row 12 (date 45123.99999999999)
# is now row 11. Its quantity changes from 360 to 220 (penalty/reward adjustment).
$ws1.Range("B11").Value = 220

# --- Sheet 2: "Monthly Trend" ---
$ws2 = $wb.Worksheets.Item(2)

# Delete the last row: date 45382.99999999999, qty 10
$ws2.Rows.Item(10).Delete()

# Row 5 (date 45138.99999999999) quantity changes from 1020 to 540
$ws2.Range("B5").Value = 540
